$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" header in column H, matching the style of the other header cells
$ws.Range("H1").Value = "Label"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("H1").VerticalAlignment = -4160    # xlTop
$ws.Range("H1").Borders.LineStyle = 1        # xlContinuous
$ws.Range("H1").Borders.Weight = 2           # xlThin

# Column H "Label" values: 0 for Control rows, 1 for MDD rows, for each of the two blocks
$labels = @(0,0,0,0,0,1,1,1,1,1,0,0,0,0,0,1,1,1,1,1)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}
